# Apply the "Office Theme" design (the stock default PowerPoint theme)
# to the presentation, replacing the custom "Integral" theme that the
# deck currently uses.
#
# PowerPoint's COM RGB() convention packs colour bytes as 0x00BBGGRR
# (blue in the high byte, red in the low byte), so every literal below
# is R | (G << 8) | (B << 16) for the target hex colour.
#
# Scheme slot order (1-based, matches ThemeColorScheme.Item(n)):
#   1 dk1=000000  2 lt1=FFFFFF  3 dk2=44546A   4 lt2=E7E6E6
#   5 accent1=5B9BD5  6 accent2=ED7D31  7 accent3=A5A5A5  8 accent4=FFC000
#   9 accent5=4472C4  10 accent6=70AD47  11 hlink=0563C1  12 folHlink=954F72

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

$officeThemeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
